# Applies the "Added experiments with memory barriers" edit:
#  1. "Таб 1,2,3" paragraph: drop the spell-check proofErr markers and
#     merge the two runs into a single run.
#  2. "Что было измерено..." paragraph: drop the paragraph-mark
#     language formatting (<w:pPr><w:rPr><w:lang .../></w:rPr></w:pPr>).
#  3. Picture paragraph: mark the run that holds the <w:drawing> as
#     <w:noProof/> (keeps Word from re-running spell/grammar check on it).
#  4. After "Стандарт не " add two empty paragraphs followed by a new
#     paragraph of text.

$d = $word.ActiveDocument

$wWmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------
# 1. "Таб 1,2,3" -- collapse proofErr/two-run markup into one clean run.
# ---------------------------------------------------------------------
$tabPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13) -eq "Таб 1,2,3") {
        $tabPara = $p
        break
    }
}
if ($tabPara -ne $null) {
    $xml = '<w:document ' + $wWmlNs + '><w:body><w:p><w:r><w:t>Таб 1,2,3</w:t></w:r></w:p></w:body></w:document>'
    $tabPara.Range.InsertXML($xml)
}

# ---------------------------------------------------------------------
# 2. "Что было измерено на примере древней функции clock()." -- remove
#    the paragraph-mark's <w:lang w:val="en-US"/> formatting.
# ---------------------------------------------------------------------
$clockPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Что было измерено")) {
        $clockPara = $p
        break
    }
}
if ($clockPara -ne $null) {
    $xml = '<w:document ' + $wWmlNs + '><w:body><w:p>' + `
        '<w:r><w:t xml:space="preserve">Что было измерено на примере древней функции </w:t></w:r>' + `
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>clock</w:t></w:r>' + `
        '<w:r><w:t>().</w:t></w:r>' + `
        '</w:p></w:body></w:document>'
    $clockPara.Range.InsertXML($xml)
}

# ---------------------------------------------------------------------
# 3. Picture paragraph -- flag the run holding the drawing as noProof.
# ---------------------------------------------------------------------
if ($d.InlineShapes.Count -ge 1) {
    $shape = $d.InlineShapes.Item(1)
    $shape.Range.NoProofing = 1
}

# ---------------------------------------------------------------------
# 4. Append two empty paragraphs and a new paragraph of text after the
#    "Стандарт не " paragraph (the final paragraph in the document).
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$endRange = $d.Range($lastPara.Range.End, $lastPara.Range.End)
$xml = '<w:document ' + $wWmlNs + '><w:body>' + `
    '<w:p/><w:p/>' + `
    '<w:p><w:r><w:t>Привилегированный режим не рассматриваем.</w:t></w:r></w:p>' + `
    '</w:body></w:document>'
$endRange.InsertXML($xml)
